$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-18 Thursday" "2024-01-19 Friday"

Replace-Text "16×33=" "68×13="
Replace-Text "99×49=" "49×14="
Replace-Text "55×49=" "14×74="
Replace-Text "45×47=" "78×43="
Replace-Text "95×46=" "50×31="

Replace-Text "99×38=" "68×59="
Replace-Text "92×47=" "29×37="
Replace-Text "13×92=" "31×11="
Replace-Text "14×60=" "51×48="
Replace-Text "65×55=" "14×17="

Replace-Text "65×22=" "67×30="
Replace-Text "35×28=" "47×92="
Replace-Text "92×96=" "81×92="
Replace-Text "73×35=" "12×50="
Replace-Text "75×51=" "50×32="

Replace-Text "77×27=" "59×81="
Replace-Text "79×20=" "46×87="
Replace-Text "67×23=" "71×91="
Replace-Text "26×32=" "40×36="
Replace-Text "31×64=" "13×18="

Replace-Text "82×36=" "77×30="
Replace-Text "76×63=" "74×13="
Replace-Text "91×43=" "88×89="
Replace-Text "82×51=" "93×76="
Replace-Text "81×89=" "39×90="
